$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two hyperlinks in C6/D6 ("Explore" / "Show (4)") along with
# their displayed text, per the diff which drops the <hyperlinks> block
# and empties those two cells (still keeping their style s="2").
$ws.Range("C6:D6").ClearContents()
$ws.Hyperlinks.Delete()

# New row 9: "passive income" / "passive.income.nadi.myfirstdrawermenuproject".
# Copy the formatting of row 8 (style s="1") down to row 9 first so the new
# cells pick up the same cell style used throughout column A/B.
$ws.Range("A8:B8").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)
$ws.Range("A9").Value = "passive income"
$ws.Range("B9").Value = "passive.income.nadi.myfirstdrawermenuproject"
$ws.Rows.Item(9).RowHeight = 24

# Widen column B to fit the longer app id text.
$ws.Columns.Item(2).ColumnWidth = 52.57

# Match the saved selection state.
$null = $ws.Range("B8").Select()
